$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5836138129234314
$ws.Range("B1").Value = 0.8064462542533875
$ws.Range("C1").Value = 7.038562297821045
$ws.Range("D1").Value = 1.81789231300354
$ws.Range("E1").Value = 1.141706585884094
